$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.196.02'
$ws.Range('E2').Value = '  -2.66%  '
$ws.Range('D3').Value = '3.010.99'
$ws.Range('E3').Value = '  -5.10%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''564.32'
$ws.Range('E5').Value = '  -4.33%  '
$ws.Range('D6').Value = '''129.06'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.008.79'
$ws.Range('E8').Value = '  -5.07%  '
$ws.Range('D9').Value = '''0.498'
$ws.Range('E9').Value = '  -2.61%  '
$ws.Range('E10').Value = '  -5.24%  '
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').Value = '''0.433'
$ws.Range('E12').Value = '  -5.45%  '
$ws.Range('D13').Value = '''0.0000224'
$ws.Range('E13').Value = '  -4.72%  '
$ws.Range('D14').Value = '''33.09'
$ws.Range('E14').Value = '  -5.18%  '
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').Value = '3.504.24'
$ws.Range('E16').Value = '  -5.18%  '
$ws.Range('D17').Value = '61.133.71'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').Value = '3.001.90'
$ws.Range('E18').Value = '  -5.39%  '
$ws.Range('D19').Value = '''6.25'
$ws.Range('E19').Value = '  -5.21%  '
$ws.Range('D20').Value = '''438.83'
$ws.Range('E20').Value = '  -4.17%  '
$ws.Range('D21').Value = '''13.23'
$ws.Range('E21').Value = '  -5.15%  '
$ws.Range('D22').Value = '''0.666'
$ws.Range('E22').Value = '  -6.35%  '
$ws.Range('D23').Value = '''7.17'
$ws.Range('E23').Value = '  -6.12%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '''12.67'
$ws.Range('E24').Value = '  -5.59%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''79.19'
$ws.Range('E25').Value = '  -4.90%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '''0.997'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').Value = '''2.50'
$ws.Range('E28').Value = '  -6.82%  '
$ws.Range('D29').Value = '''7.25'
$ws.Range('E29').Value = '  -6.80%  '
$ws.Range('D30').Value = '''1.90'
$ws.Range('E30').Value = '  -6.47%  '
$ws.Range('D31').Value = '''6.21'
$ws.Range('E31').Value = '  -9.35%  '
$ws.Range('D32').Value = '''25.64'
$ws.Range('E32').Value = '  -6.21%  '
$ws.Range('D33').Value = '''0.0944'
$ws.Range('E33').Value = '  -8.99%  '
$ws.Range('D34').Value = '''2.28'
$ws.Range('E34').Value = '  -4.01%  '
$ws.Range('D35').Value = '''0.957'
$ws.Range('E35').Value = '  -7.36%  '
$ws.Range('D36').Value = '''5.58'
$ws.Range('E36').Value = '  -4.30%  '
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('D38').Value = '0.0₃0673'
$ws.Range('E38').Value = '  -5.09%  '
$ws.Range('D39').Value = '''0.0362'
$ws.Range('E39').Value = '  -6.72%  '
$ws.Range('D40').Value = '''7.78'
$ws.Range('E40').Value = '  -3.88%  '
$ws.Range('D41').Value = '''0.109'
$ws.Range('E41').Value = '  -2.55%  '
$ws.Range('D42').Value = '''377.46'
$ws.Range('E42').Value = '  -5.86%  '
$ws.Range('D43').Value = '2.677.94'
$ws.Range('E43').Value = '  -3.46%  '
$ws.Range('D44').Value = '''2.48'
$ws.Range('E44').Value = '  -8.35%  '
$ws.Range('D46').Value = '''0.237'
$ws.Range('E46').Value = '  -6.35%  '
$ws.Range('D47').Value = '''34.27'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').Value = '''1.99'
$ws.Range('E48').Value = '  -6.78%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''118.67'
$ws.Range('E49').Value = '  -5.40%  '
$ws.Range('E50').Value = '  -3.88%  '
$ws.Range('D51').Value = '''23.52'
$ws.Range('E51').Value = '  -7.95%  '
